$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.AddShape("curvedRightArrow", 100, 100, 200, 200)
$shp.TextFrame.TextRange.Font.Color.ObjectThemeColor = 1
